$d = $word.ActiveDocument
$count = 0

if ($d.Content.Find.Execute("74-60=", $true, $true, $false, $false, $false, $true, 1, $false, "39+58=", 2)) { $count++ } else { Write-Output "FAILED: 74-60= -> 39+58=" }
if ($d.Content.Find.Execute("66-29=", $true, $true, $false, $false, $false, $true, 1, $false, "90-63=", 2)) { $count++ } else { Write-Output "FAILED: 66-29= -> 90-63=" }
if ($d.Content.Find.Execute("94-41=", $true, $true, $false, $false, $false, $true, 1, $false, "93-79=", 2)) { $count++ } else { Write-Output "FAILED: 94-41= -> 93-79=" }
if ($d.Content.Find.Execute("75+16=", $true, $true, $false, $false, $false, $true, 1, $false, "96-63=", 2)) { $count++ } else { Write-Output "FAILED: 75+16= -> 96-63=" }
if ($d.Content.Find.Execute("84-51=", $true, $true, $false, $false, $false, $true, 1, $false, "74-64=", 2)) { $count++ } else { Write-Output "FAILED: 84-51= -> 74-64=" }
if ($d.Content.Find.Execute("19+32=", $true, $true, $false, $false, $false, $true, 1, $false, "16+39=", 2)) { $count++ } else { Write-Output "FAILED: 19+32= -> 16+39=" }
if ($d.Content.Find.Execute("4+87=", $true, $true, $false, $false, $false, $true, 1, $false, "60+14=", 2)) { $count++ } else { Write-Output "FAILED: 4+87= -> 60+14=" }
if ($d.Content.Find.Execute("32+32=", $true, $true, $false, $false, $false, $true, 1, $false, "57-57=", 2)) { $count++ } else { Write-Output "FAILED: 32+32= -> 57-57=" }
if ($d.Content.Find.Execute("30-5=", $true, $true, $false, $false, $false, $true, 1, $false, "5+36=", 2)) { $count++ } else { Write-Output "FAILED: 30-5= -> 5+36=" }
if ($d.Content.Find.Execute("58+39=", $true, $true, $false, $false, $false, $true, 1, $false, "75-19=", 2)) { $count++ } else { Write-Output "FAILED: 58+39= -> 75-19=" }
if ($d.Content.Find.Execute("57+39=", $true, $true, $false, $false, $false, $true, 1, $false, "34+35=", 2)) { $count++ } else { Write-Output "FAILED: 57+39= -> 34+35=" }
if ($d.Content.Find.Execute("2+23=", $true, $true, $false, $false, $false, $true, 1, $false, "86+5=", 2)) { $count++ } else { Write-Output "FAILED: 2+23= -> 86+5=" }
if ($d.Content.Find.Execute("31-27=", $true, $true, $false, $false, $false, $true, 1, $false, "80-23=", 2)) { $count++ } else { Write-Output "FAILED: 31-27= -> 80-23=" }
if ($d.Content.Find.Execute("17+28=", $true, $true, $false, $false, $false, $true, 1, $false, "8+37=", 2)) { $count++ } else { Write-Output "FAILED: 17+28= -> 8+37=" }
if ($d.Content.Find.Execute("4+37=", $true, $true, $false, $false, $false, $true, 1, $false, "64-62=", 2)) { $count++ } else { Write-Output "FAILED: 4+37= -> 64-62=" }
if ($d.Content.Find.Execute("87-1=", $true, $true, $false, $false, $false, $true, 1, $false, "92-1=", 2)) { $count++ } else { Write-Output "FAILED: 87-1= -> 92-1=" }
if ($d.Content.Find.Execute("53-8=", $true, $true, $false, $false, $false, $true, 1, $false, "64-61=", 2)) { $count++ } else { Write-Output "FAILED: 53-8= -> 64-61=" }
if ($d.Content.Find.Execute("81+3=", $true, $true, $false, $false, $false, $true, 1, $false, "62+11=", 2)) { $count++ } else { Write-Output "FAILED: 81+3= -> 62+11=" }
if ($d.Content.Find.Execute("37-0=", $true, $true, $false, $false, $false, $true, 1, $false, "9+19=", 2)) { $count++ } else { Write-Output "FAILED: 37-0= -> 9+19=" }
if ($d.Content.Find.Execute("8+67=", $true, $true, $false, $false, $false, $true, 1, $false, "18+25=", 2)) { $count++ } else { Write-Output "FAILED: 8+67= -> 18+25=" }
if ($d.Content.Find.Execute("53+8=", $true, $true, $false, $false, $false, $true, 1, $false, "9+69=", 2)) { $count++ } else { Write-Output "FAILED: 53+8= -> 9+69=" }
if ($d.Content.Find.Execute("63-61=", $true, $true, $false, $false, $false, $true, 1, $false, "13-12=", 2)) { $count++ } else { Write-Output "FAILED: 63-61= -> 13-12=" }
if ($d.Content.Find.Execute("42-23=", $true, $true, $false, $false, $false, $true, 1, $false, "23-14=", 2)) { $count++ } else { Write-Output "FAILED: 42-23= -> 23-14=" }
if ($d.Content.Find.Execute("57+26=", $true, $true, $false, $false, $false, $true, 1, $false, "61-36=", 2)) { $count++ } else { Write-Output "FAILED: 57+26= -> 61-36=" }
if ($d.Content.Find.Execute("15+25=", $true, $true, $false, $false, $false, $true, 1, $false, "60+19=", 2)) { $count++ } else { Write-Output "FAILED: 15+25= -> 60+19=" }
if ($d.Content.Find.Execute("7+47=", $true, $true, $false, $false, $false, $true, 1, $false, "9+21=", 2)) { $count++ } else { Write-Output "FAILED: 7+47= -> 9+21=" }
if ($d.Content.Find.Execute("26+48=", $true, $true, $false, $false, $false, $true, 1, $false, "17-6=", 2)) { $count++ } else { Write-Output "FAILED: 26+48= -> 17-6=" }
if ($d.Content.Find.Execute("89-43=", $true, $true, $false, $false, $false, $true, 1, $false, "98-45=", 2)) { $count++ } else { Write-Output "FAILED: 89-43= -> 98-45=" }
if ($d.Content.Find.Execute("34+40=", $true, $true, $false, $false, $false, $true, 1, $false, "60+5=", 2)) { $count++ } else { Write-Output "FAILED: 34+40= -> 60+5=" }
if ($d.Content.Find.Execute("39+37=", $true, $true, $false, $false, $false, $true, 1, $false, "99-61=", 2)) { $count++ } else { Write-Output "FAILED: 39+37= -> 99-61=" }
if ($d.Content.Find.Execute("9+87=", $true, $true, $false, $false, $false, $true, 1, $false, "34+21=", 2)) { $count++ } else { Write-Output "FAILED: 9+87= -> 34+21=" }
if ($d.Content.Find.Execute("77-8=", $true, $true, $false, $false, $false, $true, 1, $false, "43+49=", 2)) { $count++ } else { Write-Output "FAILED: 77-8= -> 43+49=" }
if ($d.Content.Find.Execute("2+50=", $true, $true, $false, $false, $false, $true, 1, $false, "6+82=", 2)) { $count++ } else { Write-Output "FAILED: 2+50= -> 6+82=" }
if ($d.Content.Find.Execute("16+44=", $true, $true, $false, $false, $false, $true, 1, $false, "40-9=", 2)) { $count++ } else { Write-Output "FAILED: 16+44= -> 40-9=" }
if ($d.Content.Find.Execute("32+23=", $true, $true, $false, $false, $false, $true, 1, $false, "47+29=", 2)) { $count++ } else { Write-Output "FAILED: 32+23= -> 47+29=" }
if ($d.Content.Find.Execute("79-8=", $true, $true, $false, $false, $false, $true, 1, $false, "17-17=", 2)) { $count++ } else { Write-Output "FAILED: 79-8= -> 17-17=" }
if ($d.Content.Find.Execute("79-18=", $true, $true, $false, $false, $false, $true, 1, $false, "98-10=", 2)) { $count++ } else { Write-Output "FAILED: 79-18= -> 98-10=" }
if ($d.Content.Find.Execute("89-3=", $true, $true, $false, $false, $false, $true, 1, $false, "69-49=", 2)) { $count++ } else { Write-Output "FAILED: 89-3= -> 69-49=" }
if ($d.Content.Find.Execute("74+13=", $true, $true, $false, $false, $false, $true, 1, $false, "10+24=", 2)) { $count++ } else { Write-Output "FAILED: 74+13= -> 10+24=" }
if ($d.Content.Find.Execute("2+78=", $true, $true, $false, $false, $false, $true, 1, $false, "64-12=", 2)) { $count++ } else { Write-Output "FAILED: 2+78= -> 64-12=" }
if ($d.Content.Find.Execute("62+27=", $true, $true, $false, $false, $false, $true, 1, $false, "95-57=", 2)) { $count++ } else { Write-Output "FAILED: 62+27= -> 95-57=" }
if ($d.Content.Find.Execute("13+52=", $true, $true, $false, $false, $false, $true, 1, $false, "60+4=", 2)) { $count++ } else { Write-Output "FAILED: 13+52= -> 60+4=" }
if ($d.Content.Find.Execute("49-6=", $true, $true, $false, $false, $false, $true, 1, $false, "14+53=", 2)) { $count++ } else { Write-Output "FAILED: 49-6= -> 14+53=" }
if ($d.Content.Find.Execute("46+20=", $true, $true, $false, $false, $false, $true, 1, $false, "98-10=", 2)) { $count++ } else { Write-Output "FAILED: 46+20= -> 98-10=" }
if ($d.Content.Find.Execute("92+4=", $true, $true, $false, $false, $false, $true, 1, $false, "24+6=", 2)) { $count++ } else { Write-Output "FAILED: 92+4= -> 24+6=" }
if ($d.Content.Find.Execute("65+10=", $true, $true, $false, $false, $false, $true, 1, $false, "7+19=", 2)) { $count++ } else { Write-Output "FAILED: 65+10= -> 7+19=" }
if ($d.Content.Find.Execute("83-8=", $true, $true, $false, $false, $false, $true, 1, $false, "36-30=", 2)) { $count++ } else { Write-Output "FAILED: 83-8= -> 36-30=" }
if ($d.Content.Find.Execute("69-67=", $true, $true, $false, $false, $false, $true, 1, $false, "94-40=", 2)) { $count++ } else { Write-Output "FAILED: 69-67= -> 94-40=" }
if ($d.Content.Find.Execute("8+10=", $true, $true, $false, $false, $false, $true, 1, $false, "76+20=", 2)) { $count++ } else { Write-Output "FAILED: 8+10= -> 76+20=" }
if ($d.Content.Find.Execute("83-35=", $true, $true, $false, $false, $false, $true, 1, $false, "4+36=", 2)) { $count++ } else { Write-Output "FAILED: 83-35= -> 4+36=" }
if ($d.Content.Find.Execute("47+52=", $true, $true, $false, $false, $false, $true, 1, $false, "19+18=", 2)) { $count++ } else { Write-Output "FAILED: 47+52= -> 19+18=" }
if ($d.Content.Find.Execute("90-13=", $true, $true, $false, $false, $false, $true, 1, $false, "68+24=", 2)) { $count++ } else { Write-Output "FAILED: 90-13= -> 68+24=" }
if ($d.Content.Find.Execute("35+28=", $true, $true, $false, $false, $false, $true, 1, $false, "19+35=", 2)) { $count++ } else { Write-Output "FAILED: 35+28= -> 19+35=" }
if ($d.Content.Find.Execute("88-63=", $true, $true, $false, $false, $false, $true, 1, $false, "6+92=", 2)) { $count++ } else { Write-Output "FAILED: 88-63= -> 6+92=" }
if ($d.Content.Find.Execute("88-77=", $true, $true, $false, $false, $false, $true, 1, $false, "54-21=", 2)) { $count++ } else { Write-Output "FAILED: 88-77= -> 54-21=" }
if ($d.Content.Find.Execute("4+16=", $true, $true, $false, $false, $false, $true, 1, $false, "49+13=", 2)) { $count++ } else { Write-Output "FAILED: 4+16= -> 49+13=" }
if ($d.Content.Find.Execute("44+34=", $true, $true, $false, $false, $false, $true, 1, $false, "36-34=", 2)) { $count++ } else { Write-Output "FAILED: 44+34= -> 36-34=" }
if ($d.Content.Find.Execute("80-35=", $true, $true, $false, $false, $false, $true, 1, $false, "80-25=", 2)) { $count++ } else { Write-Output "FAILED: 80-35= -> 80-25=" }
if ($d.Content.Find.Execute("49-5=", $true, $true, $false, $false, $false, $true, 1, $false, "61-45=", 2)) { $count++ } else { Write-Output "FAILED: 49-5= -> 61-45=" }
if ($d.Content.Find.Execute("68-5=", $true, $true, $false, $false, $false, $true, 1, $false, "63-21=", 2)) { $count++ } else { Write-Output "FAILED: 68-5= -> 63-21=" }
if ($d.Content.Find.Execute("48+4=", $true, $true, $false, $false, $false, $true, 1, $false, "69-44=", 2)) { $count++ } else { Write-Output "FAILED: 48+4= -> 69-44=" }
if ($d.Content.Find.Execute("77-7=", $true, $true, $false, $false, $false, $true, 1, $false, "84+4=", 2)) { $count++ } else { Write-Output "FAILED: 77-7= -> 84+4=" }
if ($d.Content.Find.Execute("75-48=", $true, $true, $false, $false, $false, $true, 1, $false, "42+25=", 2)) { $count++ } else { Write-Output "FAILED: 75-48= -> 42+25=" }
if ($d.Content.Find.Execute("6+26=", $true, $true, $false, $false, $false, $true, 1, $false, "3+69=", 2)) { $count++ } else { Write-Output "FAILED: 6+26= -> 3+69=" }
if ($d.Content.Find.Execute("49-7=", $true, $true, $false, $false, $false, $true, 1, $false, "71-18=", 2)) { $count++ } else { Write-Output "FAILED: 49-7= -> 71-18=" }
if ($d.Content.Find.Execute("46+43=", $true, $true, $false, $false, $false, $true, 1, $false, "3+21=", 2)) { $count++ } else { Write-Output "FAILED: 46+43= -> 3+21=" }
if ($d.Content.Find.Execute("36+50=", $true, $true, $false, $false, $false, $true, 1, $false, "20+12=", 2)) { $count++ } else { Write-Output "FAILED: 36+50= -> 20+12=" }
if ($d.Content.Find.Execute("68+4=", $true, $true, $false, $false, $false, $true, 1, $false, "63-28=", 2)) { $count++ } else { Write-Output "FAILED: 68+4= -> 63-28=" }
if ($d.Content.Find.Execute("66-47=", $true, $true, $false, $false, $false, $true, 1, $false, "80+12=", 2)) { $count++ } else { Write-Output "FAILED: 66-47= -> 80+12=" }
if ($d.Content.Find.Execute("40-0=", $true, $true, $false, $false, $false, $true, 1, $false, "42+48=", 2)) { $count++ } else { Write-Output "FAILED: 40-0= -> 42+48=" }
if ($d.Content.Find.Execute("92-84=", $true, $true, $false, $false, $false, $true, 1, $false, "87-86=", 2)) { $count++ } else { Write-Output "FAILED: 92-84= -> 87-86=" }
if ($d.Content.Find.Execute("42-19=", $true, $true, $false, $false, $false, $true, 1, $false, "19+35=", 2)) { $count++ } else { Write-Output "FAILED: 42-19= -> 19+35=" }
if ($d.Content.Find.Execute("0+45=", $true, $true, $false, $false, $false, $true, 1, $false, "51-18=", 2)) { $count++ } else { Write-Output "FAILED: 0+45= -> 51-18=" }
if ($d.Content.Find.Execute("20+60=", $true, $true, $false, $false, $false, $true, 1, $false, "11+5=", 2)) { $count++ } else { Write-Output "FAILED: 20+60= -> 11+5=" }
if ($d.Content.Find.Execute("62-10=", $true, $true, $false, $false, $false, $true, 1, $false, "5+46=", 2)) { $count++ } else { Write-Output "FAILED: 62-10= -> 5+46=" }
if ($d.Content.Find.Execute("66+27=", $true, $true, $false, $false, $false, $true, 1, $false, "97-7=", 2)) { $count++ } else { Write-Output "FAILED: 66+27= -> 97-7=" }
if ($d.Content.Find.Execute("51+5=", $true, $true, $false, $false, $false, $true, 1, $false, "53+3=", 2)) { $count++ } else { Write-Output "FAILED: 51+5= -> 53+3=" }
if ($d.Content.Find.Execute("38-27=", $true, $true, $false, $false, $false, $true, 1, $false, "27+3=", 2)) { $count++ } else { Write-Output "FAILED: 38-27= -> 27+3=" }
if ($d.Content.Find.Execute("24+34=", $true, $true, $false, $false, $false, $true, 1, $false, "93-92=", 2)) { $count++ } else { Write-Output "FAILED: 24+34= -> 93-92=" }
if ($d.Content.Find.Execute("59-45=", $true, $true, $false, $false, $false, $true, 1, $false, "11+21=", 2)) { $count++ } else { Write-Output "FAILED: 59-45= -> 11+21=" }
if ($d.Content.Find.Execute("31-18=", $true, $true, $false, $false, $false, $true, 1, $false, "6+13=", 2)) { $count++ } else { Write-Output "FAILED: 31-18= -> 6+13=" }
if ($d.Content.Find.Execute("58+40=", $true, $true, $false, $false, $false, $true, 1, $false, "26-14=", 2)) { $count++ } else { Write-Output "FAILED: 58+40= -> 26-14=" }
if ($d.Content.Find.Execute("78-7=", $true, $true, $false, $false, $false, $true, 1, $false, "24+46=", 2)) { $count++ } else { Write-Output "FAILED: 78-7= -> 24+46=" }
if ($d.Content.Find.Execute("35-30=", $true, $true, $false, $false, $false, $true, 1, $false, "23+58=", 2)) { $count++ } else { Write-Output "FAILED: 35-30= -> 23+58=" }
if ($d.Content.Find.Execute("57+41=", $true, $true, $false, $false, $false, $true, 1, $false, "17+27=", 2)) { $count++ } else { Write-Output "FAILED: 57+41= -> 17+27=" }
if ($d.Content.Find.Execute("13+31=", $true, $true, $false, $false, $false, $true, 1, $false, "94-77=", 2)) { $count++ } else { Write-Output "FAILED: 13+31= -> 94-77=" }
if ($d.Content.Find.Execute("39+24=", $true, $true, $false, $false, $false, $true, 1, $false, "28+10=", 2)) { $count++ } else { Write-Output "FAILED: 39+24= -> 28+10=" }
if ($d.Content.Find.Execute("4+72=", $true, $true, $false, $false, $false, $true, 1, $false, "13+5=", 2)) { $count++ } else { Write-Output "FAILED: 4+72= -> 13+5=" }
if ($d.Content.Find.Execute("82-78=", $true, $true, $false, $false, $false, $true, 1, $false, "62-43=", 2)) { $count++ } else { Write-Output "FAILED: 82-78= -> 62-43=" }
if ($d.Content.Find.Execute("4+42=", $true, $true, $false, $false, $false, $true, 1, $false, "92-57=", 2)) { $count++ } else { Write-Output "FAILED: 4+42= -> 92-57=" }
if ($d.Content.Find.Execute("78+0=", $true, $true, $false, $false, $false, $true, 1, $false, "33-32=", 2)) { $count++ } else { Write-Output "FAILED: 78+0= -> 33-32=" }
if ($d.Content.Find.Execute("95-43=", $true, $true, $false, $false, $false, $true, 1, $false, "77+18=", 2)) { $count++ } else { Write-Output "FAILED: 95-43= -> 77+18=" }
if ($d.Content.Find.Execute("63+28=", $true, $true, $false, $false, $false, $true, 1, $false, "0+75=", 2)) { $count++ } else { Write-Output "FAILED: 63+28= -> 0+75=" }
if ($d.Content.Find.Execute("9+16=", $true, $true, $false, $false, $false, $true, 1, $false, "28+54=", 2)) { $count++ } else { Write-Output "FAILED: 9+16= -> 28+54=" }
if ($d.Content.Find.Execute("71+9=", $true, $true, $false, $false, $false, $true, 1, $false, "10+58=", 2)) { $count++ } else { Write-Output "FAILED: 71+9= -> 10+58=" }
if ($d.Content.Find.Execute("3+7=", $true, $true, $false, $false, $false, $true, 1, $false, "30-3=", 2)) { $count++ } else { Write-Output "FAILED: 3+7= -> 30-3=" }
if ($d.Content.Find.Execute("37+38=", $true, $true, $false, $false, $false, $true, 1, $false, "42+53=", 2)) { $count++ } else { Write-Output "FAILED: 37+38= -> 42+53=" }
if ($d.Content.Find.Execute("37+48=", $true, $true, $false, $false, $false, $true, 1, $false, "35-25=", 2)) { $count++ } else { Write-Output "FAILED: 37+48= -> 35-25=" }
if ($d.Content.Find.Execute("72-28=", $true, $true, $false, $false, $false, $true, 1, $false, "92-82=", 2)) { $count++ } else { Write-Output "FAILED: 72-28= -> 92-82=" }
if ($d.Content.Find.Execute("92-88=", $true, $true, $false, $false, $false, $true, 1, $false, "48-41=", 2)) { $count++ } else { Write-Output "FAILED: 92-88= -> 48-41=" }

Write-Output "Done applying replacements: $count / 100"
